$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-09-17 21:00:11"

for ($row = 2; $row -le 63; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
